$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values
$ws.Range("B2").Value = 11554763.3588294
$ws.Range("C2").Value = -268663126.363068
$ws.Range("D2").Value = 150283.802560066
$ws.Range("E2").Value = 113384.311737736
$ws.Range("F2").Value = 277990070.426937
$ws.Range("G2").Value = 278103454.738674
$ws.Range("H2").Value = 6323.68866933204
$ws.Range("I2").Value = 0.0407166101333371
$ws.Range("J2").Value = 49.68896867558
$ws.Range("K2").Value = 4285.68356997017
$ws.Range("L2").Value = 18049.6889686756
$ws.Range("M2").Value = 19326.0655474582
$ws.Range("N2").Value = 227425363.977265
$ws.Range("O2").Value = 210342861.472896
$ws.Range("P2").Value = 1359713894870.56
$ws.Range("Q2").Value = 1360151663096.01

# Row 3 updated values
$ws.Range("B3").Value = 4037258.72853742
$ws.Range("C3").Value = -153382718.802206
$ws.Range("D3").Value = 103217.965641601
$ws.Range("E3").Value = 78337.2153925029
$ws.Range("F3").Value = 156660162.851105
$ws.Range("G3").Value = 156738500.066497
$ws.Range("H3").Value = 6323.68866933204
$ws.Range("I3").Value = 0.0407166101333371
$ws.Range("J3").Value = 49.68896867558
$ws.Range("K3").Value = 4285.68356997017
$ws.Range("L3").Value = 18049.6889686756
$ws.Range("M3").Value = 19326.0655474582
$ws.Range("N3").Value = 73512143.862045
$ws.Range("O3").Value = 80823710.0594101
$ws.Range("P3").Value = 766261190591.106
$ws.Range("Q3").Value = 766415526445.027
